$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.325.89"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.552.82"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'210.00"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").Value = "'0.480"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'23.88"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.773.97"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "1.553.53"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "28.326.73"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "'3.62"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "'60.82"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "'227.79"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "'3.92"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "'8.93"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").Value = "'14.76"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("D30").Value = "'0.0469"
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("E31").Value = "  -4.61%  "
$ws.Range("D32").Value = "'3.17"
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").Value = "1.389.10"
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("E35").Value = "  +2.42%  "
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").Value = "'2.58"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").Value = "'1.91"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'0.777"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").Value = "'0.0460"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").Value = "'5.37"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("D46").Value = "'61.96"
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("D47").Value = "1.686.17"
$ws.Range("D48").Value = "'0.889"
$ws.Range("E48").Value = "  -7.55%  "
$ws.Range("D49").Value = "'85.64"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "'42.50"
$ws.Range("E50").Value = "  +6.93%  "
$ws.Range("E51").Value = "  +0.29%  "
